$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

# Column C ("Förändrad") holds a date serial that is bumped by one day
# for every data row (rows 2..lastRow) on every automatic update run.
$ws.Range("C2:C$lastRow").Value = 46082
